$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A (shifts former A..I to B..J) to make room for 'id'
$ws.Columns("A:A").Insert()

# Newly introduced header cells (A, K, L, M) need the same bold/
# bordered/centered style as the pre-existing header cells. Copy the
# formatting (not the value) from an existing styled header cell (B1)
# so the resulting style index matches exactly.
foreach ($colIdx in 1,11,12,13) {
    $ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,$colIdx))
}

# ---- Header row (row 1) ----
$ws.Cells.Item(1,1).Value = 'id'
$ws.Cells.Item(1,2).Value = 'Title'
$ws.Cells.Item(1,3).Value = 'Company'
$ws.Cells.Item(1,4).Value = 'Location'
$ws.Cells.Item(1,5).Value = 'Responsibilities'
$ws.Cells.Item(1,6).Value = 'Company_Summary'
$ws.Cells.Item(1,7).Value = 'Start_Date'
$ws.Cells.Item(1,8).Value = 'Deadline'
$ws.Cells.Item(1,9).Value = 'Posting_Date'
$ws.Cells.Item(1,10).Value = 'URL'
$ws.Cells.Item(1,11).Value = 'Notes'
$ws.Cells.Item(1,12).Value = 'Embedding'
$ws.Cells.Item(1,13).Value = 'Status'

# ---- Data rows ----
# Row 2
$ws.Cells.Item(2,1).Value = 3
$ws.Cells.Item(2,2).Value = '2025 Full-Time Analyst Program - EMEA.'
$ws.Cells.Item(2,3).Value = 'BlackRock.'
$ws.Cells.Item(2,4).Value = 'EMEA.'
$ws.Cells.Item(2,5).Value = 'The program is a two-year experience to empower and support Analysts in connecting their personal passions and strengths to BlackRock''s mission, principles, and purpose. The program includes an orientation, joining teams, ongoing training, and professional development to contribute to BlackRock''s collective purpose. Candidates can apply for up to two functions within the program.'
$ws.Cells.Item(2,6).Value = 'Not Found.'
$ws.Cells.Item(2,7).Value = 'Not Found.'
$ws.Cells.Item(2,8).Value = 'Not Found.'
$ws.Cells.Item(2,9).Value = '2024-07-10.'
$ws.Cells.Item(2,10).Value = 'https://blackrock.tal.net/vx/brand-3/spa-1/candidate/so/pm/1/pl/1/opp/8160-2025-Full-Time-Analyst-Program-EMEA/en-GB'
$ws.Cells.Item(2,11).Value = 'do this  ASAP'
$ws.Cells.Item(2,12).Value = '[-0.033743176609277725, -0.048345617949962616, -0.04582605138421059, 0.013140936382114887, -0.00043513879063539207, -0.015164581127464771, -0.10559264570474625, 0.007391585037112236, -0.0623038113117218, -0.08784788846969604, -0.015139295719563961, -0.04571225866675377, -0.041405774652957916, -0.028423428535461426, -0.05413978546857834, 0.036913663148880005, 0.03060143068432808, -0.12113075703382492, 0.04774594306945801, -0.08333242684602737, -0.05349944904446602, -0.026873230934143066, -0.04815082252025604, -0.02028927579522133, 0.0014181418810039759, 0.044124800711870193, -0.04913049936294556, -0.028951171785593033, -0.03853727877140045, -0.053328972309827805, 0.029082005843520164, -5.544238956645131e-05, 0.09907113760709763, 0.038419876247644424, 0.07208051532506943, 0.09523702412843704, -0.054625652730464935, 0.047536879777908325, -0.020030422136187553, -0.022603219375014305, -0.059872303158044815, -0.044038835912942886, -0.052203256636857986, 0.05772116780281067, 0.012114507146179676, -0.104239322245121, 0.08666917681694031, -0.05901888385415077, -0.037115879356861115, 0.030139373615384102, -0.07813785970211029, -0.11146333813667297, 0.015209292992949486, -0.04517224431037903, 0.021418193355202675, 0.05772911384701729, 0.03642810508608818, 0.0023434890899807215, -0.0045218993909657, -0.0284180399030447, 0.03661283850669861, -0.006058202590793371, -0.010248485952615738, 0.007699958980083466, 0.08357939124107361, -0.01504268404096365, 0.0028167208656668663, 0.03713802993297577, -0.03549402952194214, -0.06463000178337097, 0.03591427579522133, -0.09128601104021072, -0.08588709682226181, 0.044679779559373856, 0.010283307172358036, 0.13696134090423584, 0.08258771896362305, 0.01805289089679718, 0.13758234679698944, -0.0871838852763176, 0.025577416643500328, -0.055596817284822464, -0.0366859957575798, 0.018358996137976646, -0.02246084250509739, -0.01606466807425022, -0.017656702548265457, 0.05017930269241333, 0.04534928873181343, 0.0017415721667930484, -0.00966441910713911, 0.0027401060797274113, 0.018950436264276505, -0.037208981812000275, 0.016445331275463104, -0.0031954955775290728, -0.07462989538908005, -0.03491680324077606, -0.03750348836183548, 0.03597944602370262, -0.07508028298616409, -0.0034927166998386383, 0.025378765538334846, -0.02279968000948429, -0.1519404947757721, -0.0334007628262043, 0.007413189392536879, 0.043899357318878174, 0.004648895468562841, -0.002632536692544818, 0.015588360838592052, 0.05052582919597626, -0.012563342228531837, -0.0027032264042645693, 0.06694953143596649, 0.04260191321372986, -0.08911026269197464, 0.05600409582257271, 0.11271732300519943, 0.03410647064447403, 0.008458830416202545, 0.06804075837135315, 0.03352232277393341, -0.025033406913280487, -0.04733284190297127, 0.007784607354551554, -0.08863060176372528, -4.4783429091230715e-33, -0.03788352012634277, 0.011514340527355671, -0.004670565016567707, 0.0705268383026123, 0.010261401534080505, -0.02506116032600403, -0.008492971770465374, 0.04498014971613884, -0.03323822095990181, 0.08934672176837921, 0.00696595897898078, 0.11406432092189789, 0.039476796984672546, 0.04462333396077156, -0.011903868988156319, -0.04863271862268448, 0.05833371728658676, 0.0484379343688488, -0.012043274939060211, -0.02039942517876625, 0.06520663946866989, -0.10908326506614685, -0.043450139462947845, -0.005750551354140043, 0.010391389019787312, -0.03272391855716705, 0.023627661168575287, -0.0029343105852603912, 0.040900494903326035, 0.014601020142436028, -0.016780786216259003, 0.022982247173786163, -0.04299526289105415, -0.04924461245536804, 0.08587411791086197, -0.0007898451876826584, -0.025183329358696938, -0.02298629842698574, 0.024111054837703705, 0.00526385335251689, 0.002692724112421274, 0.006909716874361038, 0.013217111118137836, -0.11253891885280609, 0.03522714972496033, -0.02284696139395237, 0.02689361199736595, -0.07132551819086075, 0.0672745481133461, 0.033705100417137146, -0.07296884804964066, -0.0721363052725792, 0.01342732459306717, -0.03353992477059364, 0.04739994555711746, -0.008325611241161823, -0.01680983416736126, -0.04026823118329048, 0.04629528895020485, -0.041612330824136734, 0.02802114561200142, 0.060919277369976044, -0.11214371025562286, -0.02939079888164997, -0.11196152120828629, 0.0359991081058979, -0.0062381853349506855, 0.059845078736543655, 0.050060875713825226, 0.055870719254016876, 0.011318741366267204, -0.047046687453985214, 0.057631466537714005, 0.01550787128508091, 0.0022498241160064936, -0.0398053377866745, -0.014821785502135754, 0.03136797621846199, 0.06087091565132141, 0.04179196059703827, 0.022251516580581665, 0.03920619562268257, -0.006583237089216709, -0.03721705451607704, 0.07984662055969238, -0.0080028735101223, 0.00670284079387784, -0.0023283464834094048, 0.010731954127550125, 0.0709858387708664, -0.06200527027249336, 0.02147955633699894, -0.09665896743535995, 0.18524962663650513, 0.06818737834692001, 7.034860277688696e-34, 0.025537416338920593, -0.026821468025445938, 0.0637495368719101, -0.053762584924697876, 0.09994800388813019, 0.015113556757569313, 0.0728374570608139, 0.011677684262394905, 0.05884880572557449, 0.07701443135738373, 0.041575681418180466, 0.01486940123140812, -0.006551679689437151, -0.014281027019023895, 0.037320319563150406, -0.01384839229285717, -0.00990449357777834, -0.037720438092947006, -0.042969003319740295, 0.002496367320418358, 0.028796613216400146, -0.031100351363420486, -0.05621267482638359, 0.04142694175243378, 0.03923017159104347, 0.011183420196175575, 0.09440254420042038, -0.029290786013007164, -0.028506342321634293, 0.006825007498264313, -0.03477098047733307, -0.05119973421096802, -0.10830988734960556, 0.064421646296978, -0.06805102527141571, 0.009789342060685158, 0.015475102700293064, -0.009641352109611034, -0.04000893607735634, -0.03761162608861923, 0.06423433870077133, -0.004357772879302502, -0.01734253391623497, 0.05525343120098114, -0.043159592896699905, 0.07549973577260971, 0.05397279933094978, 0.006339446175843477, 0.032464854419231415, -0.0929417610168457, 0.023142889142036438, -0.043590448796749115, 0.041109636425971985, -0.07213156670331955, 0.02409651316702366, -0.02065020613372326, 0.08219419419765472, 0.005277693271636963, -0.03720439597964287, 0.022947115823626518, 0.03517968952655792, 0.01891339011490345, 0.02497456967830658, -0.013141949661076069, -0.04688163474202156, -0.02513211779296398, 0.030022207647562027, 0.04452364146709442, -0.0891190841794014, -0.010366219095885754, 0.03552917391061783, 0.015844298526644707, -0.053384579718112946, -0.033503204584121704, -0.11453506350517273, 0.06686800718307495, -0.05195852369070053, -0.07015703618526459, -0.029312388971447945, -0.08458513766527176, 0.013147525489330292, -0.02241380326449871, 0.02020859904587269, 0.06118416041135788, 0.022399306297302246, 0.016329802572727203, -0.028193257749080658, 0.040509942919015884, 0.051756974309682846, -0.037277016788721085, -0.09342068433761597, -0.0015448156045749784, -0.018360916525125504, 0.051383353769779205, 0.035407982766628265, -5.958699844654802e-08, -0.015901610255241394, 0.008411332964897156, 0.07155180722475052, -0.005119251552969217, 0.08424974232912064, -0.015896299853920937, 0.012719310820102692, -0.050868794322013855, 0.11915598064661026, 0.030075710266828537, -0.01595495641231537, -0.03893915191292763, -0.02242128551006317, 0.007325269281864166, -0.012719761580228806, 0.01198244746774435, 0.015816057100892067, 0.011701284907758236, 0.056125156581401825, 0.018610572442412376, 0.00427672453224659, 0.007073940709233284, 0.03820556029677391, -0.027791349217295647, -0.021298691630363464, -0.08704107254743576, -0.06978296488523483, 0.07244765013456345, 0.02060604654252529, -0.03819208964705467, 0.03742515295743942, 0.02318674698472023, 0.07983800023794174, -0.09362082928419113, 0.0017527777235955, -0.013721000403165817, 0.0550810806453228, 0.02210722491145134, -0.0848151370882988, -0.05268773436546326, -0.05256568640470505, 0.010413920506834984, 0.08267369866371155, -0.0034930824767798185, -0.05437001213431358, 0.005121666472405195, -0.13286632299423218, 0.040459927171468735, -0.020198224112391472, -0.07554101943969727, 0.02813813090324402, -0.0055683208629488945, -0.02049284055829048, 0.04921826720237732, 0.05880393460392952, 0.052497342228889465, -0.013709556311368942, -0.05425423011183739, -0.04925579950213432, 0.0165861789137125, 0.07876241952180862, -0.13480374217033386, -0.035543303936719894, -0.010141716338694096]'
$ws.Cells.Item(2,13).Value = 'assessments'

# Row 3
$ws.Cells.Item(3,1).Value = 4
$ws.Cells.Item(3,2).Value = 'Analyst'
$ws.Cells.Item(3,3).Value = 'Gain Theory'
$ws.Cells.Item(3,4).Value = 'Not Found'
$ws.Cells.Item(3,5).Value = 'As an Analyst, you will collate, analyze, and validate key performance data from clients. You will also support and engage in client calls, collaborate with the team to establish econometric models, and participate in industry-leading training academy.'
$ws.Cells.Item(3,6).Value = 'Gain Theory is a leading global marketing effectiveness and foresight consultancy that combines quality data, proprietary technology, and advanced analytics to give clients the confidence to make better informed investment decisions that drive growth.'
$ws.Cells.Item(3,7).Value = 'Not Found'
$ws.Cells.Item(3,8).Value = 'Not Found'
$ws.Cells.Item(3,9).Value = 'Not Found'
$ws.Cells.Item(3,10).Value = 'https://boards.greenhouse.io/gaintheory/jobs/7427938002'
$ws.Cells.Item(3,11).Value = ''
$ws.Cells.Item(3,12).Value = '[-0.03617524728178978, -0.01684778928756714, -0.10282022505998611, 0.01091993972659111, -0.030735304579138756, 0.03704344108700752, 0.060514915734529495, 0.031208988279104233, -0.025658220052719116, -0.02244722843170166, -0.014958664774894714, 0.03325890004634857, 0.042117640376091, 0.016822515055537224, -0.005939063150435686, 0.05527690052986145, 0.04875655472278595, -0.027031250298023224, 0.04760025069117546, -0.054056599736213684, -0.10578636825084686, -0.00741778127849102, 0.0006438311538659036, 0.0038523529656231403, 0.03493031486868858, -0.07867499440908432, -0.0044821626506745815, 0.00930619053542614, 0.01004085224121809, -0.017391767352819443, -0.06818144768476486, 0.008683830499649048, 0.06341862678527832, -0.030293060466647148, -0.028015773743391037, 0.0392865315079689, -0.022466914728283882, 0.007174388505518436, 0.03377131372690201, 0.06070803478360176, 0.008789228275418282, -0.0173904150724411, 0.05810488015413284, 0.012133161537349224, 0.006526822689920664, -0.07950325310230255, 0.037651605904102325, 0.018003756180405617, -0.050845105201005936, -0.07548589259386063, -0.1237550601363182, -0.04816557466983795, 0.0021112936083227396, -0.008025047369301319, -0.005250103306025267, 0.032785575836896896, 0.025298770517110825, 0.0014029425801709294, -0.03644296154379845, -0.061432238668203354, 0.010071387514472008, -0.017039146274328232, -0.012769605033099651, 0.02562638744711876, 0.06787535548210144, -0.028722060844302177, -0.03862445428967476, 0.053175538778305054, -0.046982817351818085, -0.04528726637363434, 0.03302427753806114, -0.07588844746351242, -0.08013007789850235, -0.009208624251186848, 0.03150466829538345, 0.026284074410796165, 0.047979313880205154, 0.07363492995500565, 0.06446101516485214, 0.04101540148258209, -0.02026582695543766, 0.04349489137530327, -0.014452176168560982, 0.006222696043550968, -0.0546291321516037, -0.02817470207810402, 0.037468623369932175, -0.01107393391430378, 0.00041245276224799454, 0.023408787325024605, 0.01687082089483738, 0.019579170271754265, -0.05793518200516701, 0.05726981163024902, -0.03346298262476921, 0.035640791058540344, -0.08764880150556564, -0.11421243846416473, 0.004476706497371197, 0.038645993918180466, -0.009535488672554493, 0.04834902659058571, -0.015075348317623138, 0.004149189684540033, -0.0675940066576004, -0.04427681863307953, -0.014226622879505157, 0.10235485434532166, 0.04566509276628494, 0.07857022434473038, -0.036807723343372345, -0.006161827594041824, -0.0020064383279532194, 0.00971301645040512, 0.05077414587140083, 0.0635719895362854, -0.06196865066885948, 0.0315643809735775, 0.003923310432583094, -0.01363882701843977, 0.06560343503952026, 0.017980076372623444, 0.05291755869984627, 0.005439455155283213, -0.07535050064325333, -0.009754320606589317, -0.09103557467460632, 5.391449035878948e-35, -0.06272536516189575, 0.04929840937256813, 0.036447882652282715, 0.08601848036050797, -0.049245964735746384, -0.024329133331775665, -0.010745538398623466, -0.016153648495674133, -0.01903172954916954, -0.04253075644373894, -0.07561866194009781, 0.09270872920751572, 0.044196657836437225, 0.048440780490636826, -0.043521322309970856, 0.012593727558851242, 0.028523104265332222, 0.11126390099525452, 0.03781392425298691, 0.007964713498950005, 0.09006647765636444, -0.037546053528785706, 0.014869979582726955, 0.016680877655744553, 0.05363954231142998, 0.04480181261897087, 0.002823798917233944, 0.07806959003210068, -0.049127258360385895, 0.003030987223610282, 0.05515991523861885, -0.0076226224191486835, -0.10653024911880493, -0.0945419892668724, -0.003395572304725647, -0.04975854605436325, 0.0011111489729955792, -0.06537288427352905, 0.03524620458483696, 0.030168600380420685, -0.04256515949964523, -0.02845599688589573, 0.011186378076672554, -0.053081247955560684, -0.02802911587059498, 0.020319504663348198, -0.0037522844504565, -0.0694185271859169, 0.0285651795566082, 0.049586955457925797, -0.098735012114048, -0.08995892852544785, 0.001711799530312419, -0.06053822487592697, -0.04883328452706337, 0.09504895657300949, -0.03074093721807003, -0.025042252615094185, 0.015500570647418499, -0.04647331312298775, -0.058911558240652084, -0.012837911956012249, -0.04124242067337036, -0.0705714076757431, -0.13450181484222412, -0.00583950849249959, 0.01450711116194725, -0.055657342076301575, 0.06628462672233582, 0.02017817832529545, 0.06089208647608757, -0.024792823940515518, 0.09216300398111343, 0.006703800987452269, 0.04512733593583107, -0.10398469865322113, -0.09723209589719772, 0.04704909399151802, 0.13979563117027283, 0.042965926229953766, 0.038985829800367355, -0.00030816736398264766, 0.03365778177976608, -0.01200962532311678, 0.041115645319223404, -0.031536512076854706, -0.0630580261349678, -0.017402181401848793, 0.08240123838186264, 0.08455193787813187, -0.09838796406984329, 0.08772924542427063, -0.1221916452050209, 0.0743587538599968, -0.023787446320056915, -2.6468174670633505e-33, 0.08179772645235062, -0.00450586574152112, -0.03758813813328743, -0.06272352486848831, 0.08349829167127609, -0.010253124870359898, -0.06175142526626587, 0.021309001371264458, 0.04320826753973961, -0.014581833966076374, 6.187948019942269e-05, -0.011732243001461029, -0.01351251546293497, 0.011899363249540329, -0.05321893468499184, -0.06559944897890091, 0.028349008411169052, -0.07214423269033432, -0.049135904759168625, -0.008081645704805851, 0.061002328991889954, -0.07125905901193619, -0.09134452790021896, -0.033068422228097916, 0.03390980884432793, -0.022683601826429367, 0.0002908932510763407, 0.0682118833065033, -0.07300160080194473, 0.0700254887342453, -0.019642017781734467, -0.02838226594030857, -0.10978225618600845, 0.03948788717389107, -0.052356742322444916, 0.07936782389879227, 0.04169036075472832, -0.047248244285583496, -0.0246355589479208, 0.057456132024526596, 0.0031747729517519474, 0.012205636128783226, 0.05930301919579506, -0.013845160603523254, 0.03357319533824921, 0.03796421363949776, 0.14784358441829681, 0.03400079533457756, 0.012254981324076653, -0.021779589354991913, 0.07487177103757858, 0.05069710686802864, 0.02819165401160717, 0.0005649125669151545, -0.049259062856435776, 0.09808260947465897, 0.092165507376194, -0.03248254209756851, -0.007856854237616062, 0.05299977958202362, 0.005591945257037878, 0.035072192549705505, -0.014052877202630043, -0.009730148129165173, -0.04300090670585632, -0.021803513169288635, 0.0006927902577444911, -0.020352711901068687, 0.020970413461327553, -0.09146761149168015, 0.05151217803359032, -0.03282065689563751, -0.010086435824632645, 0.002927085617557168, -0.04253143444657326, 0.03808117285370827, -0.0631907507777214, -0.0476609542965889, -0.013005577959120274, -0.03210096061229706, -0.02612932398915291, 0.0038683314342051744, 0.07854735851287842, 0.033021893352270126, 0.006405374966561794, 0.02149629220366478, -0.015825236216187477, -0.03443046659231186, -0.0007924208184704185, -0.0887855514883995, -0.1314254254102707, -0.016300030052661896, -0.04489191994071007, 0.02793152444064617, 0.014892063103616238, -5.684142223572053e-08, -0.100779689848423, -0.007249295245856047, 0.09053586423397064, 0.08172129839658737, -0.018719423562288284, -0.012814375571906567, -0.01683923788368702, -0.00756689952686429, 0.06095162034034729, 0.0020825499668717384, -0.04600759223103523, -0.007471881806850433, -0.07573869079351425, 0.06384991109371185, 0.05653376132249832, -0.02985851652920246, 0.023296084254980087, 0.053671807050704956, -0.043284837156534195, -0.04111054912209511, 0.11802650988101959, 0.038691721856594086, 0.05816294252872467, -0.0744866281747818, 0.043267276138067245, -0.08713839203119278, -0.07818238437175751, 0.06777147948741913, -0.01733628660440445, -0.07505732029676437, 0.019575469195842743, 0.03124314919114113, 0.05031273886561394, -0.008783383294939995, 0.05907372012734413, 0.04401947930455208, 0.02895377390086651, 0.015249097719788551, -0.01524287648499012, 0.04904765263199806, -0.04978084936738014, 0.07729530334472656, 0.07265207916498184, 0.015861356630921364, -0.055862464010715485, 0.0693650022149086, -0.1166776567697525, 0.037552010267972946, -0.031720083206892014, -0.014487185515463352, 0.06921673566102982, -0.03203704580664635, -0.06875330209732056, 0.006764883175492287, 0.02214249223470688, 0.011230901814997196, 0.03002835623919964, -0.04654210805892944, -0.03339070826768875, 0.062055692076683044, 0.013338204473257065, -0.05976705253124237, 0.009026599116623402, 0.020755186676979065]'
$ws.Cells.Item(3,13).Value = 'applied'

# Row 4
$ws.Cells.Item(4,1).Value = 5
$ws.Cells.Item(4,2).Value = 'Graduate Consultant Intake 2024.'
$ws.Cells.Item(4,3).Value = 'Quick Release (an Alten Company).'
$ws.Cells.Item(4,4).Value = 'London, England, United Kingdom.'
$ws.Cells.Item(4,5).Value = 'Work with clients to develop insights into process, systems, and data; build relationships and advise clients; shape solutions and support project plans; assist in improving consultancy division.'
$ws.Cells.Item(4,6).Value = 'Quick Release specializes in product data management and consultancy, serving various industries including automotive and aerospace, with a focus on sustainable practices.'
$ws.Cells.Item(4,7).Value = 'Not Found.'
$ws.Cells.Item(4,8).Value = 'Not Found.'
$ws.Cells.Item(4,9).Value = 'Not Found.'
$ws.Cells.Item(4,10).Value = 'https://apply.workable.com/quickrelease/j/0FF460B38D/'
$ws.Cells.Item(4,11).Value = ''
$ws.Cells.Item(4,12).Value = '[0.04312320053577423, -0.02441742829978466, 0.02473611757159233, -0.04014221578836441, -0.004613113589584827, 0.007317930925637484, 0.0013927657855674624, -0.011410627514123917, -0.06267105787992477, -0.04160746932029724, 0.03767788037657738, 0.022916579619050026, -0.04463094100356102, 0.00851225946098566, -0.012179264798760414, -0.0037409779615700245, 0.03320774435997009, -0.09867992997169495, -0.0005652212421409786, -0.027974005788564682, -0.07952183485031128, 0.0007919144118204713, -0.0297216959297657, 0.045065466314554214, -0.029927639290690422, 0.010751998052001, 0.01854005828499794, -0.0074796658009290695, 0.012399168685078621, -0.04493426904082298, 0.015951113775372505, 0.08489913493394852, -0.023065144196152687, -0.022229423746466637, 0.051801640540361404, 0.13696590065956116, -0.030268196016550064, 0.04223543033003807, 0.013052407652139664, 0.006582748610526323, -0.0344683974981308, -0.11409737914800644, -0.06624039262533188, 0.04427564516663551, 0.06575196981430054, -0.016405746340751648, 0.014030569233000278, 0.0026411954313516617, -0.05050620064139366, -0.015747586265206337, -0.11647447943687439, -0.05079912766814232, 0.007246521767228842, -0.09053082764148712, -0.05592940002679825, 0.03904743865132332, -0.02433457225561142, 0.03011462837457657, -0.05515897646546364, -0.017556309700012207, -0.023330742493271828, -0.0356125608086586, -0.021197231486439705, -0.020878512412309647, 0.04162821173667908, 0.007667718920856714, 0.02795707993209362, 0.03636974096298218, -0.04792076349258423, -0.09139443933963776, -0.023409133777022362, -0.09107837080955505, -0.07854247838258743, 0.05832641199231148, 0.017984755337238312, 0.045162100344896317, 0.02043725922703743, 0.004717602394521236, 0.06450293958187103, -0.07803269475698471, 0.08425817638635635, 0.08413401246070862, -0.0388433113694191, 0.07040176540613174, -0.08747522532939911, -0.06833934038877487, 0.026026399806141853, 0.02539004758000374, 0.006567691452801228, -0.024104658514261246, 0.0021314516197890043, -0.0503113716840744, -0.006491658743470907, 0.06189354881644249, -0.09399809688329697, 0.018462106585502625, 0.009590043686330318, -0.04428578540682793, -0.043938200920820236, 0.00048229729873128235, -0.03452175855636597, 0.03604536131024361, -0.04099122807383537, 0.03623379394412041, -0.14993703365325928, -0.044436607509851456, -0.04150080680847168, 0.09848978370428085, 0.08485368639230728, 0.01669211871922016, -0.006495494395494461, 0.06416567414999008, -0.12611441314220428, -0.015486974269151688, 0.040821436792612076, -0.07275787740945816, -0.14113055169582367, 0.002647154964506626, 0.05939466133713722, 0.04377642273902893, -0.06274452060461044, 0.1141902282834053, 0.027423866093158722, -0.020619390532374382, -0.047516874969005585, 0.024411333724856377, -0.05212864652276039, -3.0775125360618902e-33, -0.028344273567199707, 0.03494307026267052, 0.03215491399168968, 0.08425966650247574, -0.028228217735886574, -0.05256505683064461, -0.02209959551692009, 0.08558598905801773, 0.0020086478907614946, -0.0008591270307078958, 0.04468179866671562, 0.04767604172229767, -0.013792370446026325, -0.03263666108250618, 0.018465055152773857, -0.01247711107134819, 0.0480479970574379, 0.11962898820638657, -0.04131714254617691, 0.0009646588587202132, -0.08264760673046112, -0.07340347021818161, 0.02335263416171074, 0.042281556874513626, 0.10219986736774445, 0.0006158718024380505, 0.03959245607256889, 0.08345631510019302, 0.08656913042068481, 0.030615447089076042, 0.01280041690915823, -0.014686552807688713, -0.03518012538552284, -0.04027124121785164, -0.01984396204352379, -0.02415081299841404, -0.10710271447896957, -0.09056967496871948, -0.0011653798865154386, 0.005724085494875908, -0.06628144532442093, 0.04784078523516655, 0.016359133645892143, -0.044943131506443024, -0.024865711107850075, -0.0035149534232914448, -0.003192795440554619, 0.015930259600281715, 0.048065658658742905, 0.00370949599891901, -0.06354566663503647, -6.026925802871119e-06, -0.0029608779586851597, 0.025986071676015854, 0.007731297053396702, 0.008746701292693615, 0.029331795871257782, -0.04812733083963394, -0.011483453214168549, -0.026329590007662773, -0.07556634396314621, 0.08332483470439911, -0.09735066443681717, -0.06446129828691483, -0.006807639729231596, -0.03823509439826012, 0.01833161897957325, -0.02447563037276268, 0.023979373276233673, 0.024745019152760506, -0.034423697739839554, 0.0405903160572052, 0.14878104627132416, 0.002129425061866641, 0.007690496277064085, 0.003360054222866893, -0.07443057000637054, 0.0026026119012385607, 0.04113808646798134, 0.07499639689922333, 0.025251183658838272, 0.06952576339244843, 0.0215024221688509, -0.020924586802721024, 0.08480818569660187, 0.0013631752226501703, 0.023212330415844917, 0.00654195062816143, 0.03119901567697525, 0.07392148673534393, -0.0896061360836029, -0.010255781933665276, -0.10217250883579254, 0.10294262319803238, -0.013516894541680813, 7.784610319274825e-35, 0.053582336753606796, -0.049075618386268616, 0.0414024256169796, -0.016293615102767944, 0.13967818021774292, 0.014203622005879879, -0.002083816798403859, -0.017928626388311386, -0.001570410095155239, 0.07071840763092041, -0.0038221741560846567, -0.023438068106770515, -0.016728339716792107, 0.01640496961772442, -0.007156104315072298, 0.0073293522000312805, -0.013517498038709164, -0.13343745470046997, 0.0113282585516572, -0.0018696708139032125, 0.04265239089727402, -0.10393978655338287, 0.046610474586486816, 0.019517626613378525, -0.0031338389962911606, -0.012139800004661083, 0.009337016381323338, -0.027538081631064415, 0.006144431885331869, 0.010464759543538094, -0.0014804948586970568, -0.015278187580406666, -0.08894632011651993, 0.01742098294198513, -0.040977366268634796, -0.011385279707610607, -0.06272076070308685, -0.029247591271996498, -0.0012193892616778612, -0.014152689836919308, -0.007403418887406588, -0.06355687975883484, -0.014069633558392525, 0.008056963793933392, 0.0070529538206756115, -0.0477907694876194, 0.06723205745220184, -0.05428415164351463, 0.09168963134288788, -0.025716392323374748, 0.035031456500291824, 0.07903825491666794, 0.020471256226301193, -0.034903012216091156, 0.00994907971471548, -0.008261345326900482, 0.08108624070882797, -0.04249493032693863, -0.011778528802096844, 0.01890997774899006, 0.05776425078511238, -0.007232497446238995, 0.031407907605171204, -0.006449750624597073, 0.007617120165377855, -0.013522578403353691, 0.0678093433380127, 0.02116544544696808, -0.01767420768737793, -0.01779302768409252, 0.05593499168753624, 0.016777517274022102, 0.010376804508268833, -0.08329816162586212, -0.019655292853713036, -0.006902060471475124, -0.02815983258187771, -0.10932659357786179, -0.014642579481005669, -0.015124247409403324, 0.07178468257188797, -0.043648611754179, 0.015058656223118305, 0.11573637276887894, -0.031638361513614655, -0.0960984081029892, -0.0025492366403341293, 0.007143187802284956, 0.04060962423682213, 0.014369376003742218, -0.05285467579960823, 0.04302285239100456, -0.03603930026292801, 0.04438247159123421, -0.010828269645571709, -5.304393368987803e-08, -0.0554424412548542, -0.04591253027319908, 0.036294132471084595, 0.015228885225951672, 0.0346713662147522, -0.007896043360233307, -0.06131962686777115, 0.08815780282020569, 0.06478160619735718, 0.03904442861676216, 0.008517157286405563, -0.08221577107906342, -0.06201160326600075, 0.06941172480583191, 0.0517238974571228, -0.026309583336114883, -0.01792844757437706, 0.08185357600450516, -0.010912716388702393, -0.01818065531551838, 0.013741712085902691, 0.04999786242842674, 0.01454977411776781, -0.0646161213517189, 0.017293328419327736, -0.08805179595947266, 0.016691815108060837, 0.0655602365732193, -0.008528552018105984, -0.08060070872306824, -0.07382239401340485, 0.011530095711350441, 0.11427077651023865, -0.00010619968816172332, 0.045477189123630524, -0.05248109996318817, 0.053896382451057434, 0.10093006491661072, -0.021405717357993126, 0.02581324428319931, -0.07720296084880829, 0.04394708573818207, 0.01133133377879858, 0.05119842290878296, -0.02569698914885521, 0.023630138486623764, -0.016851289197802544, 0.01014233473688364, -0.033980563282966614, -0.025948621332645416, 0.0033480478450655937, -0.04023803770542145, 0.07732775807380676, 0.039283137768507004, -0.01853370852768421, 0.13016974925994873, 0.01699444092810154, -0.1443551629781723, -0.048925481736660004, 0.018619276583194733, -0.05576439946889877, -0.03166110813617706, 0.021603886038064957, 0.051080141216516495]'
$ws.Cells.Item(4,13).Value = 'assessments'

# Row 5
$ws.Cells.Item(5,1).Value = 8
$ws.Cells.Item(5,2).Value = 'Graduate Business Analyst'
$ws.Cells.Item(5,3).Value = 'Sprint Reply'
$ws.Cells.Item(5,4).Value = 'Not Found'
$ws.Cells.Item(5,5).Value = 'As a Graduate Business Analyst, you''ll support the accelerated growth of the organization by helping to develop new propositions, creating use cases for intelligent automation tools, acquiring new clients and supporting business development, and delivering projects. You''ll collaborate with clients, contribute to the selection of analysis and delivery methods, and more.'
$ws.Cells.Item(5,6).Value = 'Sprint Reply is a Reply Group company that specializes in Intelligent Process Automation, creating automation platforms to solve business problems and support clients in their transformation towards the Intelligent Enterprise.'
$ws.Cells.Item(5,7).Value = 'Not Found'
$ws.Cells.Item(5,8).Value = 'Not Found'
$ws.Cells.Item(5,9).Value = 'Not Found'
$ws.Cells.Item(5,10).Value = 'https://career5.successfactors.eu/sfcareer/jobreqcareer?jobId=9996&company=C0004534147P'
$ws.Cells.Item(5,11).Value = 'no longer hiring'
$ws.Cells.Item(5,12).Value = '[-0.05981088802218437, 0.0189900491386652, 0.016492515802383423, 0.018207842484116554, -0.011496034450829029, -0.045203909277915955, 0.037243857979774475, 0.0398775115609169, -0.008095948956906796, 0.026744283735752106, -0.022465674206614494, 0.011990281753242016, 0.0016476429300382733, -0.02359849400818348, -0.019975462928414345, 0.053920138627290726, 0.0788722112774849, -0.16677546501159668, 0.0005475830985233188, -0.045850034803152084, -0.08056791126728058, 0.012652134522795677, -0.02756085805594921, 0.0033937422558665276, -0.006411571055650711, -0.0006555786821991205, -0.009712054394185543, -0.03113487921655178, -0.040229205042123795, -0.02480539120733738, -0.002101050689816475, 0.06790163367986679, 0.06724270433187485, 0.030266601592302322, 0.009173216298222542, 0.08907274156808853, 0.020420661196112633, 0.04548445716500282, 0.07454975694417953, -0.04235474392771721, -0.05584332346916199, -0.07173612713813782, -0.05854742228984833, 0.03580138087272644, 0.14968310296535492, -0.08551254868507385, -0.05818929150700569, -0.025485320016741753, -0.0430462621152401, 0.028681159019470215, -0.08019867539405823, -0.01886710710823536, -0.03173055499792099, -0.03302228823304176, -0.07089986652135849, 0.10708341747522354, 0.07129805535078049, 0.02064019814133644, 0.011927657760679722, -0.03162037953734398, 0.050943523645401, -0.03439132496714592, 0.027114413678646088, 0.0182921402156353, 0.017749816179275513, 0.017557868734002113, -0.053198203444480896, 0.015200026333332062, -0.03859221190214157, -0.07047528028488159, 0.02130703441798687, -0.08195840567350388, -0.06890594959259033, 0.10279771685600281, 0.020355476066470146, 0.014811478555202484, 0.04990389943122864, 0.03134850040078163, 0.10249204188585281, -0.01811152510344982, -0.005524137988686562, -0.020933182910084724, -0.0618700347840786, 0.04080666974186897, -0.05369837209582329, -0.046898189932107925, 0.013497785665094852, 0.03264094144105911, 0.039694953709840775, 0.0013258049730211496, 0.025651386007666588, -0.020587455481290817, -0.03186091408133507, 0.001992353703826666, -0.04189165309071541, 0.010848830454051495, -0.1028272807598114, -0.13961507380008698, 0.02115800231695175, 0.01800469122827053, -0.02899039350450039, 0.018953727558255196, 0.0204508975148201, -0.03893548250198364, -0.11117023974657059, 0.024483446031808853, -0.005876697134226561, 0.00796586275100708, 0.03226470574736595, -0.05304831638932228, 0.019340218976140022, 0.037095069885253906, 0.004094130825251341, -0.004004690796136856, 0.04449084401130676, -0.023604659363627434, -0.03943634778261185, 0.05820551514625549, 0.1416284590959549, -0.014703693799674511, 0.03654582053422928, 0.07490339130163193, -0.0741584450006485, -0.0033940055873245, -0.002241880400106311, 0.004845946095883846, -0.034739237278699875, 2.6371136693618164e-34, -0.02685299888253212, 0.03344999998807907, 0.017619233578443527, 0.019784074276685715, 0.036386508494615555, -0.03713458031415939, -0.006577903404831886, 0.005828352179378271, -0.016523577272892, 0.006796920206397772, -0.058661457151174545, 0.032847072929143906, 0.11681701987981796, -0.01942301355302334, -0.04769044741988182, -0.04467696696519852, -0.007820798084139824, 0.054387569427490234, 0.016096584498882294, -0.07683410495519638, 0.01801171898841858, -0.034381136298179626, -0.04146558418869972, -0.028877615928649902, 0.0638023316860199, 0.037517908960580826, 0.08578255772590637, -0.015868909657001495, 0.06177820265293121, 0.014359868131577969, 0.0005326957325451076, -0.007351745385676622, -0.10060305893421173, 0.03874688595533371, 0.006454522255808115, 0.0030998748261481524, -0.0477106049656868, -0.04913349077105522, -0.04454277455806732, 0.010715478099882603, -0.04383954405784607, 0.0792321190237999, -0.021275155246257782, 0.0018179150065407157, -0.010903749614953995, 0.001903995405882597, 0.04025997966527939, -0.012936153449118137, 0.1118701919913292, -0.01595149375498295, -0.07541652768850327, -0.08002471923828125, 0.04080021008849144, -0.03723817691206932, 0.045078862458467484, 0.014789371751248837, -0.02189827524125576, -0.025669138878583908, 0.04852418228983879, -0.047548308968544006, -0.019689254462718964, 0.005373495630919933, -0.10559825599193573, -0.007970429956912994, -0.012416692450642586, -0.0296481940895319, 0.01152675412595272, 0.013326345942914486, 0.13155440986156464, 0.011887701228260994, 0.05180976167321205, -0.019787395372986794, 0.1631271243095398, 0.026950832456350327, -0.025875544175505638, -0.047707002609968185, -0.03294902294874191, 0.03434818610548973, 0.00825910922139883, -0.0074126236140728, -0.017391854897141457, 0.0051646968349814415, -0.0383303128182888, -0.02021273598074913, 0.11380956321954727, -0.026348777115345, 0.011229223571717739, 0.004946887958794832, -0.011308434419333935, 0.04839875176548958, -0.05869355425238609, 0.057516373693943024, -0.07825441658496857, 0.1760009080171585, -0.03113795816898346, -1.5217166508963106e-33, 0.016233915463089943, 0.04224241152405739, -0.03932321071624756, 0.0024332297034561634, 0.09073211997747421, -0.02466077357530594, 0.05909150093793869, 0.041277170181274414, 0.00014212154201231897, -0.003957243636250496, -0.02722417376935482, -0.011196024715900421, -0.01317940466105938, 0.02095276117324829, -0.0057257916778326035, -0.026687005534768105, 0.03753389045596123, -0.03923178091645241, -0.10381194949150085, 0.020680701360106468, -0.07095692306756973, -0.006734822876751423, -0.0054638683795928955, 0.017063135281205177, 0.06151703745126724, 0.015299835242331028, 0.0028569241985678673, 0.04646950215101242, -0.009474663995206356, 0.0007597152725793421, 0.02512151189148426, -0.007335321046411991, -0.05934663861989975, 0.03850822523236275, 0.010287811979651451, 0.046789877116680145, 0.042483892291784286, -0.06898017227649689, -0.042940739542245865, -0.019154483452439308, 0.08258963376283646, -0.03743286430835724, -0.012493298389017582, 0.04683926701545715, 0.031932733952999115, 0.04267372936010361, 0.08652608096599579, -0.029198937118053436, -0.02917170338332653, -0.04590491205453873, -0.006950246635824442, -0.029092948883771896, 0.044730328023433685, -0.032681647688150406, 0.02409617230296135, -0.016612201929092407, 0.08798631280660629, -0.07749223709106445, 0.000535507861059159, 0.06592382490634918, 0.057397302240133286, 0.002066022949293256, 0.1882529854774475, 0.0036793570034205914, 0.021426929160952568, 0.037039000540971756, 0.03848499804735184, 0.013486500829458237, -0.06632157415151596, -0.02871079184114933, 0.10535107553005219, 0.006233712192624807, -0.047770578414201736, -0.025270558893680573, 0.032885871827602386, 0.03744900971651077, -0.02986716665327549, -0.04997553676366806, -0.10861794650554657, -0.021176159381866455, 0.012448740191757679, -0.012867877259850502, -0.03242836892604828, 0.14435286819934845, -0.06069722771644592, -0.12091110646724701, -0.06262200325727463, 0.022005312144756317, 0.035873059183359146, 0.023945560678839684, -0.10021483153104782, 0.000926659326069057, -0.05098656564950943, 0.08671090006828308, 0.006773888599127531, -4.922815222130339e-08, -0.0667300820350647, 0.05031304806470871, 0.03321097046136856, 0.011705657467246056, 0.048981234431266785, 0.056806087493896484, -0.030685503035783768, -0.04931618273258209, 0.06952033936977386, -0.036585185676813126, -0.0822758823633194, -0.05838409811258316, -0.04570838063955307, 0.03891567885875702, 0.06215234100818634, 0.033850204199552536, -0.043266985565423965, -0.03838475048542023, 0.02447843924164772, 0.010171633213758469, 0.15076085925102234, 0.022516317665576935, -0.038372840732336044, -0.027323871850967407, 0.018995018675923347, -0.051650624722242355, -0.07094573229551315, 0.015109906904399395, -0.005772755481302738, -0.08021239191293716, -0.05681796744465828, 0.049950502812862396, 0.043592557311058044, -0.02590087801218033, 0.05085013061761856, 0.03860578313469887, 0.07938928157091141, -0.031407780945301056, 0.005897367838770151, -0.03129115700721741, -0.07218194007873535, -0.000721006013918668, -0.00478822086006403, -0.03440243378281593, -0.02250158227980137, -0.022221138700842857, -0.04229722544550896, -0.02907126396894455, -0.0297875814139843, 0.008924142457544804, -0.016638943925499916, -0.0899827629327774, -0.017433442175388336, 0.026915378868579865, 0.0055646817199885845, 0.01555449515581131, 0.028550870716571808, -0.1439105123281479, 0.005144818685948849, 0.01843913644552231, 0.060572072863578796, -0.015439893119037151, 0.027243705466389656, 0.006799852009862661]'
$ws.Cells.Item(5,13).Value = 'rejected'

# Row 6
$ws.Cells.Item(6,1).Value = 9
$ws.Cells.Item(6,2).Value = 'Product Analyst'
$ws.Cells.Item(6,3).Value = 'The Association'
$ws.Cells.Item(6,4).Value = 'London, United Kingdom'
$ws.Cells.Item(6,5).Value = 'Conduct research, document requirements, perform complex data analysis, support project management, assist leaders with reports, presentations, and business cases, engage with internal stakeholders, track performance, conduct market analysis, analyze data, prepare analysis and reports, and coach/ guide other analysts.'
$ws.Cells.Item(6,6).Value = 'The Association is a global organization transforming the accounting and finance profession, providing quality education, resources, and training.'
$ws.Cells.Item(6,7).Value = 'Not Found.'
$ws.Cells.Item(6,8).Value = 'Not Found.'
$ws.Cells.Item(6,9).Value = 'Jun 26, 2024.'
$ws.Cells.Item(6,10).Value = 'https://aicpa.referrals.selectminds.com/jobs/product-analyst-2883?src=JB-10180'
$ws.Cells.Item(6,11).Value = ''
$ws.Cells.Item(6,12).Value = '[-0.007402987219393253, -0.038597214967012405, -0.08856050670146942, -0.018324660137295723, -0.025279061868786812, 0.07181201130151749, 0.04710090532898903, -0.006698077078908682, -0.046853505074977875, -0.06013096868991852, -0.06078155338764191, -0.0432283915579319, -0.02750742807984352, 0.046655844897031784, -0.06080866605043411, 0.0904163345694542, 0.019151611253619194, -0.06632020324468613, -0.06365164369344711, -0.09271137416362762, -0.04659633710980415, 0.0016213511116802692, 0.06136178597807884, 0.04825560748577118, -0.1134277880191803, 0.0196820255368948, 0.02729107066988945, -0.055147428065538406, -0.014524047262966633, -0.08071961998939514, -0.01934710703790188, -0.031516142189502716, 0.1327376514673233, 0.018954578787088394, 0.05462637171149254, 0.08085997402667999, -0.028622761368751526, 0.014451204799115658, 0.03819305822253227, 0.01455290149897337, 0.048888418823480606, -0.04173155501484871, -0.06866808980703354, -0.006543566007167101, 0.08280670642852783, -0.0024190926924347878, 0.03544863313436508, -0.0010116607882082462, -0.034100279211997986, 0.03402332589030266, -0.07514391094446182, -0.03364833816885948, 0.008096820674836636, -0.07632151246070862, -0.04894067347049713, 0.024094780907034874, 0.011928597465157509, 0.0016229727771133184, 0.009490103460848331, -0.016933603212237358, 0.02176905795931816, 0.009848816320300102, 0.011310117319226265, 0.09310650825500488, 0.05101710557937622, 0.003383561735972762, -0.017749452963471413, 0.10342962294816971, -0.10405535250902176, -0.06386695057153702, -0.004214310087263584, -0.07901286333799362, -0.05266386643052101, -0.005330993328243494, -0.00507209450006485, 0.04158198460936546, -0.0015522005269303918, -0.028299480676651, 0.08061299473047256, 0.011288145557045937, 0.013311037793755531, 0.03260968625545502, -0.027831457555294037, 0.06648619472980499, -0.017316067591309547, -0.07038021832704544, 0.001676514744758606, -0.03724384307861328, 0.01492595300078392, 0.018699057400226593, 0.005599452648311853, -0.0005875892238691449, 0.037224750965833664, 0.014103531837463379, -0.014315887354314327, -0.006588765420019627, 0.019421353936195374, -0.0614875927567482, 0.03813324496150017, 0.026908142492175102, -0.0444442443549633, 0.004659028258174658, 0.005402283277362585, 0.02154795452952385, -0.20091776549816132, -0.05752479285001755, -0.04311883822083473, 0.06943659484386444, 0.035006437450647354, 0.0743420198559761, -0.04043247550725937, 0.06915277987718582, -0.12947224080562592, -0.007543888408690691, 0.038432758301496506, -0.012631651945412159, -0.05239863693714142, 0.029570987448096275, 0.04497608169913292, -0.038778968155384064, 0.006716000847518444, 0.07862145453691483, 0.07369343191385269, 0.035687290132045746, -0.048605483025312424, 0.058091647922992706, -0.08203914016485214, -3.820824266518871e-33, -0.07046844810247421, -0.004903830122202635, 0.017767589539289474, 0.055805522948503494, -0.09483411908149719, -0.00909220241010189, -0.03999119624495506, 0.042813464999198914, 0.00848609209060669, 0.03892717882990837, 0.0016294927336275578, 0.10755578428506851, 0.06148277968168259, 0.007529296446591616, -0.041645072400569916, 0.06257300823926926, 0.04758977144956589, 0.07435403019189835, -0.013892083428800106, -0.011213305406272411, -0.017458535730838776, -0.008650844916701317, 0.02771148458123207, 0.04060791805386543, 0.025246229022741318, 0.0035185127053409815, 0.04580756649374962, 0.020123135298490524, 0.01375710777938366, 0.020364943891763687, 0.07879415899515152, -0.026396634057164192, -0.10739966481924057, -0.05664205178618431, -0.03129405900835991, 0.023433303460478783, -0.07005424797534943, -0.06981830298900604, 0.046643976122140884, 0.0357075035572052, -0.0001948424760485068, -0.022175297141075134, 0.003734784899279475, -0.04988354071974754, 0.00366647494956851, 0.027461057528853416, -0.014521708711981773, -0.011985591612756252, 0.06061931326985359, 0.07180610299110413, -0.06243452802300453, -0.12471593916416168, -0.007950388826429844, -0.057095933705568314, 0.08568073064088821, -0.007441791705787182, -0.014099466614425182, -0.05546289682388306, 0.02926708571612835, -0.053757403045892715, -0.032802920788526535, 0.07851537317037582, -0.09651559591293335, -0.01821613498032093, -0.14730671048164368, -0.025900835171341896, 0.011450788006186485, 0.032944682985544205, 0.03281310573220253, -0.016772449016571045, -0.006143729202449322, -0.004497999791055918, 0.09844089299440384, 0.03495766595005989, -0.023981520906090736, -0.028331579640507698, -0.086407870054245, 0.04623958095908165, 0.0358288399875164, 0.08398541808128357, 0.0010520484065636992, 0.06593988835811615, 0.06316164135932922, -0.010419405065476894, 0.059999532997608185, -0.07711488753557205, -0.008200092241168022, 0.00829586386680603, -0.004850250668823719, 0.05640340968966484, -0.04726314917206764, 0.049971841275691986, -0.053557608276605606, 0.14681702852249146, 0.02169947512447834, -6.234447347937088e-34, 0.042005397379398346, -0.01436872873455286, 0.060636118054389954, -0.05336197093129158, 0.08918733149766922, -0.05310956761240959, 0.047328148037195206, -0.028387147933244705, 0.06280383467674255, 0.05888328328728676, -0.06318212300539017, -0.014732827432453632, 0.00556511664763093, 0.006466880906373262, 0.009403764270246029, -0.040589913725852966, -0.003717747051268816, -0.06305109709501266, -0.046177688986063004, -0.020021947100758553, -0.0012580044567584991, -0.05926746875047684, 0.0790971964597702, -0.060113660991191864, -0.0039659421890974045, -0.018935419619083405, 0.06391455233097076, -0.02374076470732689, 0.03468276560306549, 0.006889180280268192, -0.0142455343157053, -0.012540417723357677, -0.08877673000097275, 0.04333695024251938, -0.06688104569911957, -0.006184250582009554, -0.002054707147181034, -0.11383695155382156, 0.013684296980500221, -0.0736360251903534, 0.030685339123010635, -0.031210724264383316, 0.014558486640453339, 0.03828195482492447, 0.017540857195854187, 0.023527903482317924, 0.07092602550983429, -0.048216138035058975, 0.037525005638599396, -0.05500001087784767, -0.014967385679483414, 0.030481215566396713, 0.052354611456394196, -0.04921342059969902, -0.010359629057347775, 0.0596456341445446, 0.05657322332262993, -0.057265035808086395, -0.020763952285051346, -0.01013929396867752, 0.07442658394575119, 0.03944440186023712, 0.09777016192674637, 0.03507613763213158, -0.04677458480000496, -0.014269232749938965, 0.026340976357460022, 0.027518216520547867, -0.035928595811128616, -0.02949357032775879, 0.06103992089629173, -0.03194069489836693, -0.06945110857486725, -0.011843503452837467, -0.06852351129055023, 0.01442775409668684, -0.0856517106294632, -0.066163569688797, -0.019624851644039154, -0.014356200583279133, 0.020479973405599594, 0.04177515208721161, 0.03857695683836937, 0.014649275690317154, 0.01671529933810234, -0.006889489945024252, -0.0013637796510010958, 0.041048284620046616, 0.03976508229970932, -0.009853185154497623, -0.03864813223481178, -0.015160993672907352, -0.07219866663217545, 0.062399495393037796, 0.01305271778255701, -5.762692723010332e-08, -0.11090084165334702, -0.03770507499575615, 0.10757581889629364, 0.00854888278990984, 0.005938549991697073, -0.07162025570869446, -0.08692969381809235, -0.028433695435523987, 0.05463165044784546, 0.0418052151799202, 0.049404170364141464, -0.03383857384324074, -0.10096391290426254, 0.03828555718064308, 0.052850402891635895, -0.05640513077378273, -0.07418788224458694, 0.044549860060214996, 0.03268159180879593, 0.03284669667482376, 0.06628898531198502, 0.03905302658677101, 0.005937275942414999, -0.003408256219699979, -0.0029667785856872797, -0.10644195228815079, -0.03377935290336609, 0.03615732863545418, -0.007607762236148119, -0.12940151989459991, -0.017511405050754547, 0.044753026217222214, 0.05647532641887665, -0.020517148077487946, -0.041502855718135834, -0.024097148329019547, 0.08655241131782532, 0.03846803307533264, -0.062411680817604065, 0.005936815403401852, -0.10402479767799377, 0.021944185718894005, 0.05350865051150322, 0.022712375968694687, -0.029241757467389107, 0.051231976598501205, -0.04598679021000862, 0.032634466886520386, 0.013948213309049606, -0.011219627223908901, -0.02193177118897438, -0.05237076058983803, 0.0034937807358801365, 0.042847149074077606, -0.014833196997642517, 0.04822397232055664, -0.003869828535243869, -0.07922577857971191, -0.02628563530743122, 0.004022584296762943, 0.042452532798051834, -0.08728588372468948, 0.004704395309090614, 0.031919319182634354]'
$ws.Cells.Item(6,13).Value = 'not_applied'

Write-Host 'edit complete'
